# Fruta / hortaliza, semanal
#
# Insert a new weekly price-observation row for "Coliflor" (Feria Lagunitas
# de Puerto Montt) right before the current first data row for that date
# series (old row 385), pushing all subsequent rows down by one.
# This mirrors how the source dataset is updated each week: a brand new
# record is prepended to this sub-range and everything else keeps its
# previous values, just shifted down one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 385:409 down to 386:410, creating a new blank row 385
# (this also grows the sheet dimension from R409 to R410 automatically).
$ws.Rows.Item(385).Insert()

$r = 385
$ws.Cells.Item($r, 1).Value  = 4
$ws.Cells.Item($r, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item($r, 3).Value  = "Los Lagos"
$ws.Cells.Item($r, 4).Value  = 44826
$ws.Cells.Item($r, 5).Value  = 10
$ws.Cells.Item($r, 6).Value  = 100112008
$ws.Cells.Item($r, 7).Value  = "Coliflor"
$ws.Cells.Item($r, 8).Value  = "Sin especificar"
$ws.Cells.Item($r, 9).Value  = "Primera"
$ws.Cells.Item($r, 10).Value = 500
$ws.Cells.Item($r, 11).Value = 1600
$ws.Cells.Item($r, 12).Value = 1800
$ws.Cells.Item($r, 13).Value = 1700
$ws.Cells.Item($r, 14).Value = "`$/unidad"
$ws.Cells.Item($r, 15).Value = "Región Metropolitana"
$ws.Cells.Item($r, 16).Value = 1700
$ws.Cells.Item($r, 17).Value = 1
$ws.Cells.Item($r, 18).Value = "Hortaliza"
